$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated MTD Target(Tk), MTD Target(Kg), Achiv%(Tk), Achiv%(Kg) values per row
$updates = @(
    @{ Row = 2;  D = 6726;  G = 16; I = 7.437109723461195;  J = 7.875 }
    @{ Row = 3;  D = 16144; G = 39; I = 3.541129831516353;  J = 3.692307692307692 }
    @{ Row = 4;  D = 7700;  G = 17; I = 77.2987012987013;   J = 79.41176470588235 }
    @{ Row = 5;  D = 45030; G = 54; I = 40.73868532089718;  J = 42.65925925925927 }
    @{ Row = 6;  D = 23757; G = 28; I = 12.20145641284674;  J = 13 }
    @{ Row = 7;  D = 9629;  G = 21; I = 521.7040191089417;  J = 542.5714285714286 }
    @{ Row = 8;  D = 3212;  G = 7;  I = 638.6849315068494;  J = 664.7142857142858 }
    @{ Row = 9;  D = 3212;  G = 7;  I = 370.6102117061021;  J = 385.7142857142857 }
    @{ Row = 10; D = 41143; G = 74; I = 0.6418345769632744; J = 0.6756756756756757 }
    @{ Row = 11; D = 13574; G = 30; I = 5.009577132753794;  J = 5.333333333333334 }
    @{ Row = 12; D = 6787;  G = 15; I = 68.63120671872697;  J = 73.06666666666666 }
    @{ Row = 13; D = 17498; G = 21; I = 8.282946622471139;  J = 8.666666666666668 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Range("D$r").Value = $u.D
    $ws.Range("G$r").Value = $u.G
    $ws.Range("I$r").Value = $u.I
    $ws.Range("J$r").Value = $u.J
}
